$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New "Pass/Fail" results in column F, mirroring the D column's wrapped/top-aligned style.
$ws.Range("F2").Value = "pass"
$ws.Range("F3").Value = "fail"
$ws.Range("F4").Value = "fail"
$ws.Range("F5").Value = "pass"

$ws.Range("F2:F5").VerticalAlignment = -4160
$ws.Range("F2:F5").WrapText = $true

# Scroll the view so column C is at the left edge and F3 is the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F3").Select()
